# Auto-generated edit script: shifts D/J/K/L/M/N/P/Q values for rows 86-221 down by one,
# inserts new observation data at row 86, and appends a duplicated row 222 (old row 221 data).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 86
$ws.Cells.Item(86, 4).Value = 44771
$ws.Cells.Item(86, 10).Value = 1600
$ws.Cells.Item(86, 11).Value = 1400
$ws.Cells.Item(86, 12).Value = 1600
$ws.Cells.Item(86, 13).Value = 1500
$ws.Cells.Item(86, 16).Value = 250
# Row 87
$ws.Cells.Item(87, 4).Value = 44209
$ws.Cells.Item(87, 10).Value = 2700
$ws.Cells.Item(87, 11).Value = 1000
$ws.Cells.Item(87, 12).Value = 1100
$ws.Cells.Item(87, 13).Value = 1050
$ws.Cells.Item(87, 16).Value = 175
# Row 88
$ws.Cells.Item(88, 4).Value = 44323
$ws.Cells.Item(88, 10).Value = 2880
# Row 89
$ws.Cells.Item(89, 4).Value = 44270
$ws.Cells.Item(89, 10).Value = 2200
$ws.Cells.Item(89, 11).Value = 900
$ws.Cells.Item(89, 12).Value = 1000
$ws.Cells.Item(89, 13).Value = 950
$ws.Cells.Item(89, 14).Value = '$/paquete 6 unidades'
$ws.Cells.Item(89, 16).Value = 158
$ws.Cells.Item(89, 17).Value = 6
# Row 90
$ws.Cells.Item(90, 4).Value = 44603
$ws.Cells.Item(90, 10).Value = 880
$ws.Cells.Item(90, 11).Value = 3500
$ws.Cells.Item(90, 12).Value = 4000
$ws.Cells.Item(90, 13).Value = 3750
$ws.Cells.Item(90, 14).Value = '$/paquete 36 unidades'
$ws.Cells.Item(90, 16).Value = 104
$ws.Cells.Item(90, 17).Value = 36
# Row 91
$ws.Cells.Item(91, 4).Value = 44407
$ws.Cells.Item(91, 10).Value = 3200
$ws.Cells.Item(91, 11).Value = 900
$ws.Cells.Item(91, 12).Value = 1000
$ws.Cells.Item(91, 13).Value = 950
$ws.Cells.Item(91, 14).Value = '$/paquete 6 unidades'
$ws.Cells.Item(91, 16).Value = 158
$ws.Cells.Item(91, 17).Value = 6
# Row 92
$ws.Cells.Item(92, 4).Value = 44599
$ws.Cells.Item(92, 10).Value = 1000
$ws.Cells.Item(92, 11).Value = 3800
$ws.Cells.Item(92, 12).Value = 4000
$ws.Cells.Item(92, 13).Value = 3900
$ws.Cells.Item(92, 14).Value = '$/paquete 36 unidades'
$ws.Cells.Item(92, 16).Value = 108
$ws.Cells.Item(92, 17).Value = 36
# Row 93
$ws.Cells.Item(93, 4).Value = 44272
$ws.Cells.Item(93, 10).Value = 3100
$ws.Cells.Item(93, 11).Value = 800
$ws.Cells.Item(93, 12).Value = 1000
$ws.Cells.Item(93, 13).Value = 900
$ws.Cells.Item(93, 14).Value = '$/paquete 6 unidades'
$ws.Cells.Item(93, 16).Value = 150
$ws.Cells.Item(93, 17).Value = 6
# Row 94
$ws.Cells.Item(94, 4).Value = 44580
$ws.Cells.Item(94, 10).Value = 1600
$ws.Cells.Item(94, 11).Value = 4000
$ws.Cells.Item(94, 12).Value = 4500
$ws.Cells.Item(94, 13).Value = 4250
$ws.Cells.Item(94, 14).Value = '$/paquete 36 unidades'
$ws.Cells.Item(94, 16).Value = 118
$ws.Cells.Item(94, 17).Value = 36
# Row 95
$ws.Cells.Item(95, 4).Value = 44524
$ws.Cells.Item(95, 10).Value = 3200
$ws.Cells.Item(95, 11).Value = 900
$ws.Cells.Item(95, 12).Value = 1000
$ws.Cells.Item(95, 13).Value = 950
$ws.Cells.Item(95, 16).Value = 158
# Row 96
$ws.Cells.Item(96, 4).Value = 44699
$ws.Cells.Item(96, 11).Value = 1200
$ws.Cells.Item(96, 12).Value = 1400
$ws.Cells.Item(96, 13).Value = 1300
$ws.Cells.Item(96, 16).Value = 217
# Row 97
$ws.Cells.Item(97, 4).Value = 44706
$ws.Cells.Item(97, 10).Value = 1300
$ws.Cells.Item(97, 11).Value = 1400
$ws.Cells.Item(97, 12).Value = 1600
$ws.Cells.Item(97, 13).Value = 1500
$ws.Cells.Item(97, 16).Value = 250
# Row 98
$ws.Cells.Item(98, 4).Value = 44552
$ws.Cells.Item(98, 10).Value = 3000
# Row 99
$ws.Cells.Item(99, 4).Value = 44307
$ws.Cells.Item(99, 10).Value = 3200
# Row 100
$ws.Cells.Item(100, 4).Value = 44491
$ws.Cells.Item(100, 10).Value = 3000
$ws.Cells.Item(100, 11).Value = 900
$ws.Cells.Item(100, 12).Value = 1000
$ws.Cells.Item(100, 13).Value = 950
$ws.Cells.Item(100, 16).Value = 158
# Row 101
$ws.Cells.Item(101, 4).Value = 44232
$ws.Cells.Item(101, 10).Value = 2600
$ws.Cells.Item(101, 11).Value = 1000
$ws.Cells.Item(101, 12).Value = 1100
$ws.Cells.Item(101, 13).Value = 1050
$ws.Cells.Item(101, 14).Value = '$/paquete 6 unidades'
$ws.Cells.Item(101, 16).Value = 175
$ws.Cells.Item(101, 17).Value = 6
# Row 102
$ws.Cells.Item(102, 4).Value = 44629
$ws.Cells.Item(102, 10).Value = 2900
$ws.Cells.Item(102, 11).Value = 8000
$ws.Cells.Item(102, 12).Value = 9000
$ws.Cells.Item(102, 13).Value = 8500
$ws.Cells.Item(102, 14).Value = '$/paquete 36 unidades'
$ws.Cells.Item(102, 16).Value = 236
$ws.Cells.Item(102, 17).Value = 36
# Row 103
$ws.Cells.Item(103, 4).Value = 44225
$ws.Cells.Item(103, 10).Value = 2600
$ws.Cells.Item(103, 11).Value = 1000
$ws.Cells.Item(103, 12).Value = 1100
$ws.Cells.Item(103, 13).Value = 1050
$ws.Cells.Item(103, 16).Value = 175
# Row 104
$ws.Cells.Item(104, 4).Value = 44587
# Row 105
$ws.Cells.Item(105, 4).Value = 44249
$ws.Cells.Item(105, 10).Value = 2800
# Row 106
$ws.Cells.Item(106, 4).Value = 44445
$ws.Cells.Item(106, 10).Value = 3100
$ws.Cells.Item(106, 11).Value = 900
$ws.Cells.Item(106, 12).Value = 1000
$ws.Cells.Item(106, 13).Value = 950
$ws.Cells.Item(106, 14).Value = '$/paquete 6 unidades'
$ws.Cells.Item(106, 16).Value = 158
$ws.Cells.Item(106, 17).Value = 6
# Row 107
$ws.Cells.Item(107, 4).Value = 44573
$ws.Cells.Item(107, 10).Value = 400
$ws.Cells.Item(107, 11).Value = 4000
$ws.Cells.Item(107, 12).Value = 4500
$ws.Cells.Item(107, 13).Value = 4250
$ws.Cells.Item(107, 14).Value = '$/paquete 36 unidades'
$ws.Cells.Item(107, 16).Value = 118
$ws.Cells.Item(107, 17).Value = 36
# Row 108
$ws.Cells.Item(108, 4).Value = 44302
$ws.Cells.Item(108, 10).Value = 2800
$ws.Cells.Item(108, 11).Value = 900
$ws.Cells.Item(108, 12).Value = 1000
$ws.Cells.Item(108, 13).Value = 950
$ws.Cells.Item(108, 16).Value = 158
# Row 109
$ws.Cells.Item(109, 4).Value = 44664
$ws.Cells.Item(109, 10).Value = 1600
$ws.Cells.Item(109, 11).Value = 1100
$ws.Cells.Item(109, 12).Value = 1200
$ws.Cells.Item(109, 13).Value = 1150
$ws.Cells.Item(109, 16).Value = 192
# Row 110
$ws.Cells.Item(110, 4).Value = 44505
$ws.Cells.Item(110, 10).Value = 3100
$ws.Cells.Item(110, 11).Value = 900
$ws.Cells.Item(110, 13).Value = 950
$ws.Cells.Item(110, 16).Value = 158
# Row 111
$ws.Cells.Item(111, 4).Value = 44167
$ws.Cells.Item(111, 10).Value = 2700
$ws.Cells.Item(111, 11).Value = 800
$ws.Cells.Item(111, 13).Value = 900
$ws.Cells.Item(111, 16).Value = 150
# Row 112
$ws.Cells.Item(112, 4).Value = 44284
$ws.Cells.Item(112, 10).Value = 2400
$ws.Cells.Item(112, 11).Value = 900
$ws.Cells.Item(112, 13).Value = 950
$ws.Cells.Item(112, 16).Value = 158
# Row 113
$ws.Cells.Item(113, 4).Value = 44179
$ws.Cells.Item(113, 10).Value = 2200
$ws.Cells.Item(113, 11).Value = 800
$ws.Cells.Item(113, 13).Value = 900
$ws.Cells.Item(113, 16).Value = 150
# Row 114
$ws.Cells.Item(114, 4).Value = 44291
$ws.Cells.Item(114, 10).Value = 2600
$ws.Cells.Item(114, 11).Value = 900
$ws.Cells.Item(114, 12).Value = 1000
$ws.Cells.Item(114, 13).Value = 950
$ws.Cells.Item(114, 16).Value = 158
# Row 115
$ws.Cells.Item(115, 4).Value = 44561
$ws.Cells.Item(115, 10).Value = 3200
$ws.Cells.Item(115, 11).Value = 1000
$ws.Cells.Item(115, 12).Value = 1200
$ws.Cells.Item(115, 13).Value = 1100
$ws.Cells.Item(115, 16).Value = 183
# Row 116
$ws.Cells.Item(116, 4).Value = 44312
$ws.Cells.Item(116, 10).Value = 2700
$ws.Cells.Item(116, 11).Value = 900
$ws.Cells.Item(116, 12).Value = 1000
$ws.Cells.Item(116, 13).Value = 950
$ws.Cells.Item(116, 16).Value = 158
# Row 117
$ws.Cells.Item(117, 4).Value = 44221
$ws.Cells.Item(117, 10).Value = 2600
$ws.Cells.Item(117, 11).Value = 1000
$ws.Cells.Item(117, 12).Value = 1100
$ws.Cells.Item(117, 13).Value = 1050
$ws.Cells.Item(117, 16).Value = 175
# Row 118
$ws.Cells.Item(118, 4).Value = 44543
$ws.Cells.Item(118, 10).Value = 3000
# Row 119
$ws.Cells.Item(119, 4).Value = 44342
$ws.Cells.Item(119, 10).Value = 3300
$ws.Cells.Item(119, 11).Value = 900
$ws.Cells.Item(119, 12).Value = 1000
$ws.Cells.Item(119, 13).Value = 950
$ws.Cells.Item(119, 16).Value = 158
# Row 120
$ws.Cells.Item(120, 4).Value = 44396
$ws.Cells.Item(120, 10).Value = 2800
$ws.Cells.Item(120, 11).Value = 1000
$ws.Cells.Item(120, 12).Value = 1200
$ws.Cells.Item(120, 13).Value = 1100
$ws.Cells.Item(120, 16).Value = 183
# Row 121
$ws.Cells.Item(121, 4).Value = 44435
$ws.Cells.Item(121, 10).Value = 9720
# Row 122
$ws.Cells.Item(122, 4).Value = 44498
# Row 123
$ws.Cells.Item(123, 4).Value = 44519
$ws.Cells.Item(123, 10).Value = 3000
$ws.Cells.Item(123, 11).Value = 900
$ws.Cells.Item(123, 12).Value = 1000
$ws.Cells.Item(123, 13).Value = 950
$ws.Cells.Item(123, 14).Value = '$/paquete 6 unidades'
$ws.Cells.Item(123, 16).Value = 158
$ws.Cells.Item(123, 17).Value = 6
# Row 124
$ws.Cells.Item(124, 4).Value = 44582
$ws.Cells.Item(124, 10).Value = 1700
$ws.Cells.Item(124, 11).Value = 4000
$ws.Cells.Item(124, 12).Value = 4500
$ws.Cells.Item(124, 13).Value = 4250
$ws.Cells.Item(124, 14).Value = '$/paquete 36 unidades'
$ws.Cells.Item(124, 16).Value = 118
$ws.Cells.Item(124, 17).Value = 36
# Row 125
$ws.Cells.Item(125, 4).Value = 44274
$ws.Cells.Item(125, 10).Value = 2700
$ws.Cells.Item(125, 11).Value = 900
$ws.Cells.Item(125, 12).Value = 1000
$ws.Cells.Item(125, 13).Value = 950
$ws.Cells.Item(125, 16).Value = 158
# Row 126
$ws.Cells.Item(126, 4).Value = 44242
$ws.Cells.Item(126, 10).Value = 2800
$ws.Cells.Item(126, 11).Value = 1000
$ws.Cells.Item(126, 12).Value = 1100
$ws.Cells.Item(126, 13).Value = 1050
$ws.Cells.Item(126, 16).Value = 175
# Row 127
$ws.Cells.Item(127, 4).Value = 44692
$ws.Cells.Item(127, 10).Value = 1300
$ws.Cells.Item(127, 11).Value = 1100
$ws.Cells.Item(127, 12).Value = 1200
$ws.Cells.Item(127, 13).Value = 1150
$ws.Cells.Item(127, 16).Value = 192
# Row 128
$ws.Cells.Item(128, 4).Value = 44550
$ws.Cells.Item(128, 10).Value = 3000
$ws.Cells.Item(128, 11).Value = 900
$ws.Cells.Item(128, 13).Value = 950
$ws.Cells.Item(128, 16).Value = 158
# Row 129
$ws.Cells.Item(129, 4).Value = 44286
$ws.Cells.Item(129, 10).Value = 3200
$ws.Cells.Item(129, 11).Value = 800
$ws.Cells.Item(129, 12).Value = 1000
$ws.Cells.Item(129, 13).Value = 900
$ws.Cells.Item(129, 16).Value = 150
# Row 130
$ws.Cells.Item(130, 4).Value = 44223
$ws.Cells.Item(130, 10).Value = 2800
$ws.Cells.Item(130, 11).Value = 1000
$ws.Cells.Item(130, 12).Value = 1100
$ws.Cells.Item(130, 13).Value = 1050
$ws.Cells.Item(130, 16).Value = 175
# Row 131
$ws.Cells.Item(131, 4).Value = 44265
$ws.Cells.Item(131, 10).Value = 3200
$ws.Cells.Item(131, 11).Value = 900
$ws.Cells.Item(131, 12).Value = 1000
$ws.Cells.Item(131, 13).Value = 950
$ws.Cells.Item(131, 16).Value = 158
# Row 132
$ws.Cells.Item(132, 4).Value = 44762
$ws.Cells.Item(132, 10).Value = 1000
$ws.Cells.Item(132, 11).Value = 1400
$ws.Cells.Item(132, 12).Value = 1600
$ws.Cells.Item(132, 13).Value = 1500
$ws.Cells.Item(132, 16).Value = 250
# Row 133
$ws.Cells.Item(133, 4).Value = 44454
$ws.Cells.Item(133, 10).Value = 3360
# Row 134
$ws.Cells.Item(134, 4).Value = 44351
$ws.Cells.Item(134, 10).Value = 2960
$ws.Cells.Item(134, 11).Value = 900
$ws.Cells.Item(134, 12).Value = 1000
$ws.Cells.Item(134, 13).Value = 950
$ws.Cells.Item(134, 16).Value = 158
# Row 135
$ws.Cells.Item(135, 4).Value = 44645
$ws.Cells.Item(135, 10).Value = 2100
$ws.Cells.Item(135, 11).Value = 1100
$ws.Cells.Item(135, 12).Value = 1200
$ws.Cells.Item(135, 13).Value = 1150
$ws.Cells.Item(135, 16).Value = 192
# Row 136
$ws.Cells.Item(136, 4).Value = 44459
$ws.Cells.Item(136, 10).Value = 2800
# Row 137
$ws.Cells.Item(137, 4).Value = 44428
$ws.Cells.Item(137, 10).Value = 3120
# Row 138
$ws.Cells.Item(138, 4).Value = 44333
$ws.Cells.Item(138, 10).Value = 2760
$ws.Cells.Item(138, 11).Value = 900
$ws.Cells.Item(138, 12).Value = 1000
$ws.Cells.Item(138, 13).Value = 950
$ws.Cells.Item(138, 16).Value = 158
# Row 139
$ws.Cells.Item(139, 4).Value = 44767
$ws.Cells.Item(139, 10).Value = 1200
$ws.Cells.Item(139, 11).Value = 1400
$ws.Cells.Item(139, 12).Value = 1600
$ws.Cells.Item(139, 13).Value = 1500
$ws.Cells.Item(139, 16).Value = 250
# Row 140
$ws.Cells.Item(140, 4).Value = 44529
# Row 141
$ws.Cells.Item(141, 4).Value = 44533
$ws.Cells.Item(141, 10).Value = 3000
$ws.Cells.Item(141, 11).Value = 900
$ws.Cells.Item(141, 12).Value = 1000
$ws.Cells.Item(141, 13).Value = 950
$ws.Cells.Item(141, 16).Value = 158
# Row 142
$ws.Cells.Item(142, 4).Value = 44687
$ws.Cells.Item(142, 10).Value = 2100
$ws.Cells.Item(142, 11).Value = 1100
$ws.Cells.Item(142, 12).Value = 1200
$ws.Cells.Item(142, 13).Value = 1150
$ws.Cells.Item(142, 14).Value = '$/paquete 6 unidades'
$ws.Cells.Item(142, 16).Value = 192
$ws.Cells.Item(142, 17).Value = 6
# Row 143
$ws.Cells.Item(143, 4).Value = 44601
$ws.Cells.Item(143, 10).Value = 800
$ws.Cells.Item(143, 11).Value = 3500
$ws.Cells.Item(143, 12).Value = 4000
$ws.Cells.Item(143, 13).Value = 3750
$ws.Cells.Item(143, 14).Value = '$/paquete 36 unidades'
$ws.Cells.Item(143, 16).Value = 104
$ws.Cells.Item(143, 17).Value = 36
# Row 144
$ws.Cells.Item(144, 4).Value = 44431
$ws.Cells.Item(144, 10).Value = 3200
# Row 145
$ws.Cells.Item(145, 4).Value = 44263
$ws.Cells.Item(145, 10).Value = 2600
$ws.Cells.Item(145, 11).Value = 900
$ws.Cells.Item(145, 12).Value = 1000
$ws.Cells.Item(145, 13).Value = 950
$ws.Cells.Item(145, 16).Value = 158
# Row 146
$ws.Cells.Item(146, 4).Value = 44753
$ws.Cells.Item(146, 10).Value = 2000
$ws.Cells.Item(146, 11).Value = 1400
$ws.Cells.Item(146, 12).Value = 1600
$ws.Cells.Item(146, 13).Value = 1500
$ws.Cells.Item(146, 16).Value = 250
# Row 147
$ws.Cells.Item(147, 4).Value = 44267
$ws.Cells.Item(147, 10).Value = 2400
$ws.Cells.Item(147, 11).Value = 900
$ws.Cells.Item(147, 12).Value = 1000
$ws.Cells.Item(147, 13).Value = 950
$ws.Cells.Item(147, 16).Value = 158
# Row 148
$ws.Cells.Item(148, 4).Value = 44678
$ws.Cells.Item(148, 10).Value = 1200
$ws.Cells.Item(148, 11).Value = 1100
$ws.Cells.Item(148, 12).Value = 1200
$ws.Cells.Item(148, 13).Value = 1150
$ws.Cells.Item(148, 16).Value = 192
# Row 149
$ws.Cells.Item(149, 4).Value = 44344
$ws.Cells.Item(149, 10).Value = 2960
$ws.Cells.Item(149, 11).Value = 900
$ws.Cells.Item(149, 12).Value = 1000
$ws.Cells.Item(149, 13).Value = 950
$ws.Cells.Item(149, 16).Value = 158
# Row 150
$ws.Cells.Item(150, 4).Value = 44711
$ws.Cells.Item(150, 10).Value = 880
$ws.Cells.Item(150, 11).Value = 1400
$ws.Cells.Item(150, 12).Value = 1600
$ws.Cells.Item(150, 13).Value = 1500
$ws.Cells.Item(150, 14).Value = '$/paquete 6 unidades'
$ws.Cells.Item(150, 16).Value = 250
$ws.Cells.Item(150, 17).Value = 6
# Row 151
$ws.Cells.Item(151, 4).Value = 44594
$ws.Cells.Item(151, 10).Value = 800
$ws.Cells.Item(151, 11).Value = 4000
$ws.Cells.Item(151, 12).Value = 4500
$ws.Cells.Item(151, 13).Value = 4250
$ws.Cells.Item(151, 14).Value = '$/paquete 36 unidades'
$ws.Cells.Item(151, 16).Value = 118
$ws.Cells.Item(151, 17).Value = 36
# Row 152
$ws.Cells.Item(152, 4).Value = 44260
$ws.Cells.Item(152, 10).Value = 2700
$ws.Cells.Item(152, 11).Value = 900
$ws.Cells.Item(152, 12).Value = 1000
$ws.Cells.Item(152, 13).Value = 950
$ws.Cells.Item(152, 16).Value = 158
# Row 153
$ws.Cells.Item(153, 4).Value = 44195
$ws.Cells.Item(153, 10).Value = 2400
$ws.Cells.Item(153, 12).Value = 1100
$ws.Cells.Item(153, 13).Value = 1050
$ws.Cells.Item(153, 16).Value = 175
# Row 154
$ws.Cells.Item(154, 4).Value = 44568
$ws.Cells.Item(154, 10).Value = 3100
$ws.Cells.Item(154, 11).Value = 1000
$ws.Cells.Item(154, 13).Value = 1100
$ws.Cells.Item(154, 16).Value = 183
# Row 155
$ws.Cells.Item(155, 4).Value = 44634
$ws.Cells.Item(155, 10).Value = 1400
$ws.Cells.Item(155, 11).Value = 1100
$ws.Cells.Item(155, 12).Value = 1200
$ws.Cells.Item(155, 13).Value = 1150
$ws.Cells.Item(155, 16).Value = 192
# Row 156
$ws.Cells.Item(156, 4).Value = 44725
$ws.Cells.Item(156, 10).Value = 1200
$ws.Cells.Item(156, 11).Value = 1400
$ws.Cells.Item(156, 12).Value = 1600
$ws.Cells.Item(156, 13).Value = 1500
$ws.Cells.Item(156, 16).Value = 250
# Row 157
$ws.Cells.Item(157, 4).Value = 44365
$ws.Cells.Item(157, 10).Value = 2900
$ws.Cells.Item(157, 11).Value = 900
$ws.Cells.Item(157, 12).Value = 1000
$ws.Cells.Item(157, 13).Value = 950
$ws.Cells.Item(157, 16).Value = 158
# Row 158
$ws.Cells.Item(158, 4).Value = 44736
$ws.Cells.Item(158, 10).Value = 2400
$ws.Cells.Item(158, 11).Value = 1400
$ws.Cells.Item(158, 12).Value = 1600
$ws.Cells.Item(158, 13).Value = 1500
$ws.Cells.Item(158, 16).Value = 250
# Row 159
$ws.Cells.Item(159, 4).Value = 44648
$ws.Cells.Item(159, 10).Value = 2000
$ws.Cells.Item(159, 11).Value = 1100
$ws.Cells.Item(159, 12).Value = 1200
$ws.Cells.Item(159, 13).Value = 1150
$ws.Cells.Item(159, 16).Value = 192
# Row 160
$ws.Cells.Item(160, 4).Value = 44314
$ws.Cells.Item(160, 10).Value = 3200
$ws.Cells.Item(160, 11).Value = 900
$ws.Cells.Item(160, 12).Value = 1000
$ws.Cells.Item(160, 13).Value = 950
$ws.Cells.Item(160, 16).Value = 158
# Row 161
$ws.Cells.Item(161, 4).Value = 44676
$ws.Cells.Item(161, 10).Value = 800
$ws.Cells.Item(161, 11).Value = 1100
$ws.Cells.Item(161, 12).Value = 1200
$ws.Cells.Item(161, 13).Value = 1150
$ws.Cells.Item(161, 14).Value = '$/paquete 6 unidades'
$ws.Cells.Item(161, 16).Value = 192
$ws.Cells.Item(161, 17).Value = 6
# Row 162
$ws.Cells.Item(162, 4).Value = 44615
$ws.Cells.Item(162, 10).Value = 1000
$ws.Cells.Item(162, 11).Value = 8000
$ws.Cells.Item(162, 12).Value = 9000
$ws.Cells.Item(162, 13).Value = 8500
$ws.Cells.Item(162, 14).Value = '$/paquete 36 unidades'
$ws.Cells.Item(162, 16).Value = 236
$ws.Cells.Item(162, 17).Value = 36
# Row 163
$ws.Cells.Item(163, 4).Value = 44239
$ws.Cells.Item(163, 10).Value = 2600
$ws.Cells.Item(163, 11).Value = 1000
$ws.Cells.Item(163, 12).Value = 1100
$ws.Cells.Item(163, 13).Value = 1050
$ws.Cells.Item(163, 16).Value = 175
# Row 164
$ws.Cells.Item(164, 4).Value = 44638
$ws.Cells.Item(164, 10).Value = 1800
$ws.Cells.Item(164, 11).Value = 1100
$ws.Cells.Item(164, 12).Value = 1200
$ws.Cells.Item(164, 13).Value = 1150
$ws.Cells.Item(164, 16).Value = 192
# Row 165
$ws.Cells.Item(165, 4).Value = 44218
$ws.Cells.Item(165, 10).Value = 2600
# Row 166
$ws.Cells.Item(166, 4).Value = 44230
$ws.Cells.Item(166, 10).Value = 2800
$ws.Cells.Item(166, 12).Value = 1100
$ws.Cells.Item(166, 13).Value = 1050
$ws.Cells.Item(166, 16).Value = 175
# Row 167
$ws.Cells.Item(167, 4).Value = 44403
$ws.Cells.Item(167, 10).Value = 3200
$ws.Cells.Item(167, 11).Value = 1000
$ws.Cells.Item(167, 12).Value = 1200
$ws.Cells.Item(167, 13).Value = 1100
$ws.Cells.Item(167, 16).Value = 183
# Row 168
$ws.Cells.Item(168, 4).Value = 44176
$ws.Cells.Item(168, 10).Value = 2400
$ws.Cells.Item(168, 11).Value = 800
$ws.Cells.Item(168, 13).Value = 900
$ws.Cells.Item(168, 16).Value = 150
# Row 169
$ws.Cells.Item(169, 4).Value = 44463
# Row 170
$ws.Cells.Item(170, 4).Value = 44382
$ws.Cells.Item(170, 10).Value = 3000
$ws.Cells.Item(170, 11).Value = 900
$ws.Cells.Item(170, 12).Value = 1000
$ws.Cells.Item(170, 13).Value = 950
$ws.Cells.Item(170, 16).Value = 158
# Row 171
$ws.Cells.Item(171, 4).Value = 44685
$ws.Cells.Item(171, 10).Value = 1300
$ws.Cells.Item(171, 11).Value = 1100
$ws.Cells.Item(171, 12).Value = 1200
$ws.Cells.Item(171, 13).Value = 1150
$ws.Cells.Item(171, 16).Value = 192
# Row 172
$ws.Cells.Item(172, 4).Value = 44379
$ws.Cells.Item(172, 10).Value = 3000
$ws.Cells.Item(172, 11).Value = 900
$ws.Cells.Item(172, 12).Value = 1000
$ws.Cells.Item(172, 13).Value = 950
$ws.Cells.Item(172, 16).Value = 158
# Row 173
$ws.Cells.Item(173, 4).Value = 44673
$ws.Cells.Item(173, 10).Value = 2000
$ws.Cells.Item(173, 11).Value = 1100
$ws.Cells.Item(173, 12).Value = 1200
$ws.Cells.Item(173, 13).Value = 1150
$ws.Cells.Item(173, 16).Value = 192
# Row 174
$ws.Cells.Item(174, 4).Value = 44503
$ws.Cells.Item(174, 10).Value = 2800
$ws.Cells.Item(174, 11).Value = 900
$ws.Cells.Item(174, 12).Value = 1000
$ws.Cells.Item(174, 13).Value = 950
$ws.Cells.Item(174, 16).Value = 158
# Row 175
$ws.Cells.Item(175, 4).Value = 44746
$ws.Cells.Item(175, 10).Value = 1000
$ws.Cells.Item(175, 11).Value = 1400
$ws.Cells.Item(175, 12).Value = 1600
$ws.Cells.Item(175, 13).Value = 1500
$ws.Cells.Item(175, 16).Value = 250
# Row 176
$ws.Cells.Item(176, 4).Value = 44172
$ws.Cells.Item(176, 10).Value = 2400
$ws.Cells.Item(176, 11).Value = 800
$ws.Cells.Item(176, 13).Value = 900
$ws.Cells.Item(176, 16).Value = 150
# Row 177
$ws.Cells.Item(177, 4).Value = 44328
$ws.Cells.Item(177, 10).Value = 3240
# Row 178
$ws.Cells.Item(178, 4).Value = 44554
$ws.Cells.Item(178, 10).Value = 3000
# Row 179
$ws.Cells.Item(179, 4).Value = 44452
$ws.Cells.Item(179, 10).Value = 3100
# Row 180
$ws.Cells.Item(180, 4).Value = 44370
$ws.Cells.Item(180, 10).Value = 3300
# Row 181
$ws.Cells.Item(181, 4).Value = 44421
$ws.Cells.Item(181, 10).Value = 3100
# Row 182
$ws.Cells.Item(182, 4).Value = 44389
$ws.Cells.Item(182, 10).Value = 3000
# Row 183
$ws.Cells.Item(183, 4).Value = 44417
$ws.Cells.Item(183, 10).Value = 3300
# Row 184
$ws.Cells.Item(184, 4).Value = 44354
$ws.Cells.Item(184, 10).Value = 2600
$ws.Cells.Item(184, 11).Value = 900
$ws.Cells.Item(184, 12).Value = 1000
$ws.Cells.Item(184, 13).Value = 950
$ws.Cells.Item(184, 14).Value = '$/paquete 6 unidades'
$ws.Cells.Item(184, 16).Value = 158
$ws.Cells.Item(184, 17).Value = 6
# Row 185
$ws.Cells.Item(185, 4).Value = 44655
$ws.Cells.Item(185, 10).Value = 400
$ws.Cells.Item(185, 11).Value = 5500
$ws.Cells.Item(185, 12).Value = 6000
$ws.Cells.Item(185, 13).Value = 5750
$ws.Cells.Item(185, 14).Value = '$/paquete 36 unidades'
$ws.Cells.Item(185, 16).Value = 160
$ws.Cells.Item(185, 17).Value = 36
# Row 186
$ws.Cells.Item(186, 4).Value = 44384
$ws.Cells.Item(186, 10).Value = 3320
# Row 187
$ws.Cells.Item(187, 4).Value = 44515
$ws.Cells.Item(187, 10).Value = 3000
$ws.Cells.Item(187, 11).Value = 900
$ws.Cells.Item(187, 12).Value = 1000
$ws.Cells.Item(187, 13).Value = 950
$ws.Cells.Item(187, 16).Value = 158
# Row 188
$ws.Cells.Item(188, 4).Value = 44200
$ws.Cells.Item(188, 10).Value = 2500
$ws.Cells.Item(188, 11).Value = 1000
$ws.Cells.Item(188, 12).Value = 1100
$ws.Cells.Item(188, 13).Value = 1050
$ws.Cells.Item(188, 16).Value = 175
# Row 189
$ws.Cells.Item(189, 4).Value = 44748
$ws.Cells.Item(189, 10).Value = 1400
$ws.Cells.Item(189, 11).Value = 1400
$ws.Cells.Item(189, 12).Value = 1600
$ws.Cells.Item(189, 13).Value = 1500
$ws.Cells.Item(189, 16).Value = 250
# Row 190
$ws.Cells.Item(190, 4).Value = 44424
$ws.Cells.Item(190, 10).Value = 3120
$ws.Cells.Item(190, 11).Value = 900
$ws.Cells.Item(190, 12).Value = 1000
$ws.Cells.Item(190, 13).Value = 950
$ws.Cells.Item(190, 16).Value = 158
# Row 191
$ws.Cells.Item(191, 4).Value = 44202
$ws.Cells.Item(191, 11).Value = 1000
$ws.Cells.Item(191, 12).Value = 1100
$ws.Cells.Item(191, 13).Value = 1050
$ws.Cells.Item(191, 16).Value = 175
# Row 192
$ws.Cells.Item(192, 4).Value = 44162
$ws.Cells.Item(192, 10).Value = 2400
$ws.Cells.Item(192, 11).Value = 800
$ws.Cells.Item(192, 12).Value = 1000
$ws.Cells.Item(192, 13).Value = 900
$ws.Cells.Item(192, 16).Value = 150
# Row 193
$ws.Cells.Item(193, 4).Value = 44235
$ws.Cells.Item(193, 11).Value = 1000
$ws.Cells.Item(193, 12).Value = 1100
$ws.Cells.Item(193, 13).Value = 1050
$ws.Cells.Item(193, 16).Value = 175
# Row 194
$ws.Cells.Item(194, 4).Value = 44708
$ws.Cells.Item(194, 11).Value = 1400
$ws.Cells.Item(194, 12).Value = 1600
$ws.Cells.Item(194, 13).Value = 1500
$ws.Cells.Item(194, 16).Value = 250
# Row 195
$ws.Cells.Item(195, 4).Value = 44174
$ws.Cells.Item(195, 10).Value = 2600
$ws.Cells.Item(195, 11).Value = 800
$ws.Cells.Item(195, 13).Value = 900
$ws.Cells.Item(195, 16).Value = 150
# Row 196
$ws.Cells.Item(196, 4).Value = 44419
# Row 197
$ws.Cells.Item(197, 4).Value = 44412
$ws.Cells.Item(197, 10).Value = 3400
$ws.Cells.Item(197, 11).Value = 900
$ws.Cells.Item(197, 12).Value = 1000
$ws.Cells.Item(197, 13).Value = 950
$ws.Cells.Item(197, 16).Value = 158
# Row 198
$ws.Cells.Item(198, 4).Value = 44237
$ws.Cells.Item(198, 11).Value = 1000
$ws.Cells.Item(198, 12).Value = 1100
$ws.Cells.Item(198, 13).Value = 1050
$ws.Cells.Item(198, 16).Value = 175
# Row 199
$ws.Cells.Item(199, 4).Value = 44253
$ws.Cells.Item(199, 10).Value = 2800
# Row 200
$ws.Cells.Item(200, 4).Value = 44494
$ws.Cells.Item(200, 10).Value = 3000
# Row 201
$ws.Cells.Item(201, 4).Value = 44487
$ws.Cells.Item(201, 10).Value = 2800
# Row 202
$ws.Cells.Item(202, 4).Value = 44356
$ws.Cells.Item(202, 10).Value = 3360
$ws.Cells.Item(202, 11).Value = 900
$ws.Cells.Item(202, 12).Value = 1000
$ws.Cells.Item(202, 13).Value = 950
$ws.Cells.Item(202, 16).Value = 158
# Row 203
$ws.Cells.Item(203, 4).Value = 44204
$ws.Cells.Item(203, 10).Value = 2400
$ws.Cells.Item(203, 11).Value = 1000
$ws.Cells.Item(203, 12).Value = 1100
$ws.Cells.Item(203, 13).Value = 1050
$ws.Cells.Item(203, 16).Value = 175
# Row 204
$ws.Cells.Item(204, 4).Value = 44484
$ws.Cells.Item(204, 10).Value = 3060
# Row 205
$ws.Cells.Item(205, 4).Value = 44298
$ws.Cells.Item(205, 10).Value = 2600
$ws.Cells.Item(205, 11).Value = 900
$ws.Cells.Item(205, 12).Value = 1000
$ws.Cells.Item(205, 13).Value = 950
$ws.Cells.Item(205, 16).Value = 158
# Row 206
$ws.Cells.Item(206, 4).Value = 44769
$ws.Cells.Item(206, 10).Value = 1360
$ws.Cells.Item(206, 11).Value = 1400
$ws.Cells.Item(206, 12).Value = 1600
$ws.Cells.Item(206, 13).Value = 1500
$ws.Cells.Item(206, 14).Value = '$/paquete 6 unidades'
$ws.Cells.Item(206, 16).Value = 250
$ws.Cells.Item(206, 17).Value = 6
# Row 207
$ws.Cells.Item(207, 4).Value = 44596
$ws.Cells.Item(207, 10).Value = 800
$ws.Cells.Item(207, 11).Value = 4000
$ws.Cells.Item(207, 12).Value = 4500
$ws.Cells.Item(207, 13).Value = 4250
$ws.Cells.Item(207, 14).Value = '$/paquete 36 unidades'
$ws.Cells.Item(207, 16).Value = 118
$ws.Cells.Item(207, 17).Value = 36
# Row 208
$ws.Cells.Item(208, 4).Value = 44340
$ws.Cells.Item(208, 10).Value = 2800
# Row 209
$ws.Cells.Item(209, 4).Value = 44496
$ws.Cells.Item(209, 10).Value = 2000
# Row 210
$ws.Cells.Item(210, 4).Value = 44377
$ws.Cells.Item(210, 10).Value = 3200
# Row 211
$ws.Cells.Item(211, 4).Value = 44512
# Row 212
$ws.Cells.Item(212, 4).Value = 44557
$ws.Cells.Item(212, 11).Value = 900
$ws.Cells.Item(212, 13).Value = 950
$ws.Cells.Item(212, 16).Value = 158
# Row 213
$ws.Cells.Item(213, 4).Value = 44279
$ws.Cells.Item(213, 10).Value = 3000
$ws.Cells.Item(213, 11).Value = 800
$ws.Cells.Item(213, 13).Value = 900
$ws.Cells.Item(213, 16).Value = 150
# Row 214
$ws.Cells.Item(214, 4).Value = 44517
$ws.Cells.Item(214, 10).Value = 3100
# Row 215
$ws.Cells.Item(215, 4).Value = 44547
# Row 216
$ws.Cells.Item(216, 4).Value = 44321
$ws.Cells.Item(216, 10).Value = 3200
# Row 217
$ws.Cells.Item(217, 4).Value = 44438
$ws.Cells.Item(217, 10).Value = 3100
$ws.Cells.Item(217, 11).Value = 900
$ws.Cells.Item(217, 12).Value = 1000
$ws.Cells.Item(217, 13).Value = 950
$ws.Cells.Item(217, 14).Value = '$/paquete 6 unidades'
$ws.Cells.Item(217, 16).Value = 158
$ws.Cells.Item(217, 17).Value = 6
# Row 218
$ws.Cells.Item(218, 4).Value = 44657
$ws.Cells.Item(218, 10).Value = 1000
$ws.Cells.Item(218, 11).Value = 6500
$ws.Cells.Item(218, 12).Value = 7000
$ws.Cells.Item(218, 13).Value = 6750
$ws.Cells.Item(218, 14).Value = '$/paquete 36 unidades'
$ws.Cells.Item(218, 16).Value = 188
$ws.Cells.Item(218, 17).Value = 36
# Row 219
$ws.Cells.Item(219, 4).Value = 44391
$ws.Cells.Item(219, 10).Value = 3360
$ws.Cells.Item(219, 11).Value = 900
$ws.Cells.Item(219, 13).Value = 950
$ws.Cells.Item(219, 16).Value = 158
# Row 220
$ws.Cells.Item(220, 4).Value = 44186
$ws.Cells.Item(220, 10).Value = 2200
$ws.Cells.Item(220, 11).Value = 800
$ws.Cells.Item(220, 13).Value = 900
$ws.Cells.Item(220, 16).Value = 150
# Row 221
$ws.Cells.Item(221, 4).Value = 44489
$ws.Cells.Item(221, 10).Value = 3200

# New row 222 (duplicate of former row 221 data, since all rows shifted down by one)
$ws.Cells.Item(222, 1).Value = 8
$ws.Cells.Item(222, 2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(222, 3).Value = 'Coquimbo'
$ws.Cells.Item(222, 4).Value = 44358
$ws.Range("D222").NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(222, 5).Value = 4
$ws.Cells.Item(222, 6).Value = 100112037
$ws.Cells.Item(222, 7).Value = 'Cebollín'
$ws.Cells.Item(222, 8).Value = 'Sin especificar'
$ws.Cells.Item(222, 9).Value = 'Primera'
$ws.Cells.Item(222, 10).Value = 3000
$ws.Cells.Item(222, 11).Value = 900
$ws.Cells.Item(222, 12).Value = 1000
$ws.Cells.Item(222, 13).Value = 950
$ws.Cells.Item(222, 14).Value = '$/paquete 6 unidades'
$ws.Cells.Item(222, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(222, 16).Value = 158
$ws.Cells.Item(222, 17).Value = 6
$ws.Cells.Item(222, 18).Value = 'Hortaliza'
